$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure D-column price cells retain their original text formatting
# (prices are stored as text, e.g. "613.41", not as numbers)
foreach ($addr in @('D2', 'D3', 'D5', 'D6', 'D7', 'D9', 'D10', 'D11', 'D12', 'D13', 'D14', 'D15', 'D16', 'D17', 'D19', 'D20', 'D21', 'D22', 'D23', 'D25', 'D26', 'D27', 'D29', 'D30', 'D31', 'D32', 'D33', 'D34', 'D35', 'D36', 'D37', 'D38', 'D39', 'D41', 'D42', 'D43', 'D44', 'D45', 'D46', 'D47', 'D48', 'D49', 'D51')) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range('D2').Value = '69.108.19'
$ws.Range('E2').Value = '  +2.76%  '
$ws.Range('D3').Value = '3.725.67'
$ws.Range('E3').Value = '  +1.35%  '
$ws.Range('E4').Value = '  -0.11%  '
$ws.Range('D5').Value = '613.41'
$ws.Range('E5').Value = '  +8.33%  '
$ws.Range('D6').Value = '192.07'
$ws.Range('E6').Value = '  +13.11%  '
$ws.Range('D7').Value = '0.641'
$ws.Range('E7').Value = '  +3.95%  '
$ws.Range('E8').Value = '  +0.06%  '
$ws.Range('D9').Value = '0.729'
$ws.Range('E9').Value = '  +4.91%  '
$ws.Range('D10').Value = '0.162'
$ws.Range('E10').Value = '  +1.01%  '
$ws.Range('D11').Value = '60.38'
$ws.Range('E11').Value = '  +17.86%  '
$ws.Range('D12').Value = '0.0000291'
$ws.Range('E12').Value = '  +0.92%  '
$ws.Range('D13').Value = '10.56'
$ws.Range('E13').Value = '  +1.71%  '
$ws.Range('D14').Value = '4.328.88'
$ws.Range('E14').Value = '  +1.16%  '
$ws.Range('D15').Value = '3.735.68'
$ws.Range('E15').Value = '  +0.51%  '
$ws.Range('D16').Value = '1.16'
$ws.Range('E16').Value = '  +4.46%  '
$ws.Range('D17').Value = '19.57'
$ws.Range('E17').Value = '  +2.69%  '
$ws.Range('E18').Value = '  +0.68%  '
$ws.Range('D19').Value = '13.01'
$ws.Range('E19').Value = '  +2.85%  '
$ws.Range('D20').Value = '69.021.89'
$ws.Range('E20').Value = '  +2.78%  '
$ws.Range('D21').Value = '413.36'
$ws.Range('E21').Value = '  +3.25%  '
$ws.Range('D22').Value = '4.59'
$ws.Range('E22').Value = '  +3.31%  '
$ws.Range('D23').Value = '90.40'
$ws.Range('E23').Value = '  +4.22%  '
$ws.Range('E24').Value = '  +3.38%  '
$ws.Range('D25').Value = '11.45'
$ws.Range('E25').Value = '  +9.79%  '
$ws.Range('D26').Value = '13.02'
$ws.Range('E26').Value = '  +4.00%  '
$ws.Range('D27').Value = '3.83'
$ws.Range('E27').Value = '  +3.03%  '
$ws.Range('E28').Value = '  +1.57%  '
$ws.Range('D29').Value = '9.83'
$ws.Range('E29').Value = '  +5.57%  '
$ws.Range('D30').Value = '33.02'
$ws.Range('E30').Value = '  +2.54%  '
$ws.Range('D31').Value = '7.87'
$ws.Range('E31').Value = '  +6.05%  '
$ws.Range('D32').Value = '12.81'
$ws.Range('E32').Value = '  +3.76%  '
$ws.Range('D33').Value = '650.00'
$ws.Range('E33').Value = '  +11.81%  '
$ws.Range('D34').Value = '0.123'
$ws.Range('E34').Value = '  +7.17%  '
$ws.Range('D35').Value = '46.20'
$ws.Range('E35').Value = '  +8.81%  '
$ws.Range('D36').Value = '67.06'
$ws.Range('E36').Value = '  +4.84%  '
$ws.Range('D37').Value = '0.0₃0838'
$ws.Range('E37').Value = '  -3.09%  '
$ws.Range('D38').Value = '0.417'
$ws.Range('E38').Value = '  +7.30%  '
$ws.Range('D39').Value = '0.999'
$ws.Range('E39').Value = '  +0.05%  '
$ws.Range('D41').Value = '0.141'
$ws.Range('E41').Value = '  +6.18%  '
$ws.Range('D42').Value = '3.09'
$ws.Range('E42').Value = '  +4.01%  '
$ws.Range('D43').Value = '0.0451'
$ws.Range('E43').Value = '  +4.50%  '
$ws.Range('D44').Value = '2.64'
$ws.Range('E44').Value = '  +5.21%  '
$ws.Range('D45').Value = '2.909.94'
$ws.Range('E45').Value = '  +7.72%  '
$ws.Range('D46').Value = '0.140'
$ws.Range('E46').Value = '  +6.19%  '
$ws.Range('D47').Value = '9.23'
$ws.Range('E47').Value = '  +1.92%  '
$ws.Range('D48').Value = '2.75'
$ws.Range('E48').Value = '  +2.90%  '
$ws.Range('D49').Value = '143.54'
$ws.Range('E49').Value = '  +0.93%  '
$ws.Range('E50').Value = '  -0.14%  '
$ws.Range('D51').Value = '2.61'
$ws.Range('E51').Value = '  -7.06%  '
